# Update vm_pu.xlsx results for the 380 kV case (B column slack voltage changed
# from 1.05 to 1.02 p.u., with recomputed bus voltage magnitudes in C:F and I:M)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B"=1.02; "C"=1.05019399380407; "D"=1.05783290532045; "E"=1.057187736487439; "F"=1.067625180481724; "I"=1.049160186023271; "J"=1.055228439264246; "K"=1.060566794428925; "L"=1.05992339209416; "M"=1.070332541542478 }
    3 = @{ "B"=1.02; "C"=1.051099140355649; "D"=1.058565776587103; "E"=1.057986622290977; "F"=1.068490187668519; "I"=1.04942122303769; "J"=1.055783386719374; "K"=1.061113785461762; "L"=1.06053610267233; "M"=1.071013247691224 }
    4 = @{ "B"=1.02; "C"=1.051685383885528; "D"=1.059040430719538; "E"=1.058504402466648; "F"=1.06905081019511; "I"=1.049589069366782; "J"=1.05614236634828; "K"=1.06146749172106; "L"=1.06093275705777; "M"=1.07145396915327 }
    5 = @{ "B"=1.02; "C"=1.051931971740138; "D"=1.059240078431376; "E"=1.058722278545592; "F"=1.069286710867086; "I"=1.049659376964506; "J"=1.056293254248414; "K"=1.061616132702162; "L"=1.06109955467821; "M"=1.071639309055421 }
    6 = @{ "B"=1.02; "C"=1.05197338257282; "D"=1.059273606160732; "E"=1.058758872654368; "F"=1.069326332183812; "I"=1.049671166944901; "J"=1.056318587365683; "K"=1.061641086791157; "L"=1.061127563282144; "M"=1.071670431952797 }
    7 = @{ "B"=1.02; "C"=1.05168867829117; "D"=1.059043098018878; "E"=1.058507312948584; "F"=1.06905396146892; "I"=1.049590009823232; "J"=1.056144382627869; "K"=1.061469478095143; "L"=1.060934985644006; "M"=1.07145644543778 }
    8 = @{ "B"=1.02; "C"=1.050499776960109; "D"=1.058080491062007; "E"=1.057457547193261; "F"=1.067917325291439; "I"=1.049248624170785; "J"=1.055416008116475; "K"=1.0607517002229; "L"=1.060130420454088; "M"=1.070562534876355 }
    9 = @{ "B"=1.02; "C"=1.04840907777859; "D"=1.056387675566138; "E"=1.055614289921758; "F"=1.065921433990181; "I"=1.04863895738673; "J"=1.054131739797114; "K"=1.059485149646746; "L"=1.058714183482511; "M"=1.068989397214391 }
    10 = @{ "B"=1.02; "C"=1.047018247647768; "D"=1.055261527213209; "E"=1.054389962166734; "F"=1.064595652126908; "I"=1.048227109447058; "J"=1.053275105981339; "K"=1.058639689324034; "L"=1.057771117143966; "M"=1.067942099414913 }
    11 = @{ "B"=1.02; "C"=1.046416723841862; "D"=1.054774480684231; "E"=1.053860904946502; "F"=1.064022736639709; "I"=1.048047502237816; "J"=1.052904080627928; "K"=1.058273351594534; "L"=1.057363033994001; "M"=1.067488972044367 }
    12 = @{ "B"=1.02; "C"=1.046193399827565; "D"=1.054593659131273; "E"=1.05366455415484; "F"=1.063810105824604; "I"=1.047980597376858; "J"=1.052766251825586; "K"=1.058137241425456; "L"=1.057211495708183; "M"=1.067320715784948 }
    13 = @{ "B"=1.02; "C"=1.046241298646179; "D"=1.054632441905351; "E"=1.053706664590611; "F"=1.063855707886791; "I"=1.047994957316939; "J"=1.052795817162524; "K"=1.05816643912956; "L"=1.057243999265268; "M"=1.067356804793127 }
    14 = @{ "B"=1.02; "C"=1.046398261583891; "D"=1.054759532089052; "E"=1.053844671161204; "F"=1.064005156920053; "I"=1.048041975748205; "J"=1.052892687924205; "K"=1.058262101413407; "L"=1.057350506928515; "M"=1.067475062790723 }
    15 = @{ "B"=1.02; "C"=1.046494986074845; "D"=1.054837848415365; "E"=1.05392972342645; "F"=1.06409726070923; "I"=1.04807092010097; "J"=1.052952371447508; "K"=1.058321037366223; "L"=1.057416135414283; "M"=1.067547932832194 }
    16 = @{ "B"=1.02; "C"=1.047058183974388; "D"=1.055293863310255; "E"=1.054425096947097; "F"=1.064633699159397; "I"=1.048239002597632; "J"=1.053299727760757; "K"=1.058663996834274; "L"=1.057798206093429; "M"=1.067972179680708 }
    17 = @{ "B"=1.02; "C"=1.047411655442685; "D"=1.055580066661255; "E"=1.054736123253007; "F"=1.064970503743199; "I"=1.048344095630779; "J"=1.053517589968736; "K"=1.058879060749372; "L"=1.058037942341203; "M"=1.068238395762396 }
    18 = @{ "B"=1.02; "C"=1.047617898251002; "D"=1.055747060390702; "E"=1.054917644190344; "F"=1.065167067659303; "I"=1.048405271616679; "J"=1.053644655847802; "K"=1.059004479899301; "L"=1.058177802480608; "M"=1.06839370968942 }
    19 = @{ "B"=1.02; "C"=1.047688233362288; "D"=1.055804010413341; "E"=1.054979555823516; "F"=1.065234109735496; "I"=1.048426110162119; "J"=1.053687980389311; "K"=1.059047240481894; "L"=1.058225495552041; "M"=1.068446673549065 }
    20 = @{ "B"=1.02; "C"=1.047373724165137; "D"=1.05554935393013; "E"=1.054702742247416; "F"=1.064934356253467; "I"=1.048332832860503; "J"=1.053494216385952; "K"=1.058855988899618; "L"=1.058012218207742; "M"=1.068209829708201 }
    21 = @{ "B"=1.02; "C"=1.046352036899408; "D"=1.054722104721984; "E"=1.0538040271085; "F"=1.063961143046133; "I"=1.048028135251941; "J"=1.052864162261856; "K"=1.058233932239894; "L"=1.057319141907561; "M"=1.067440237216682 }
    22 = @{ "B"=1.02; "C"=1.045710291182707; "D"=1.054202497267394; "E"=1.053239922461246; "F"=1.063350261544884; "I"=1.047835456524006; "J"=1.052467944884127; "K"=1.057842612697048; "L"=1.056883620766571; "M"=1.066956685223589 }
    23 = @{ "B"=1.02; "C"=1.046050432608053; "D"=1.054477901473633; "E"=1.053538874104981; "F"=1.063674004485855; "I"=1.047937703501519; "J"=1.052677994161084; "K"=1.058050077886957; "L"=1.057114475260625; "M"=1.067212994352019 }
    24 = @{ "B"=1.02; "C"=1.047390863464026; "D"=1.055563231517983; "E"=1.05471782536411; "F"=1.064950689403808; "I"=1.048337922401719; "J"=1.053504777931265; "K"=1.058866414149673; "L"=1.058023841756049; "M"=1.068222737371761 }
    25 = @{ "B"=1.02; "C"=1.048949056165238; "D"=1.056824893730155; "E"=1.056090028222427; "F"=1.066436578943906; "I"=1.048797526458342; "J"=1.054463839187442; "K"=1.059812781355031; "L"=1.059080127887999; "M"=1.069395839739695 }
}

foreach ($rowNum in $data.Keys) {
    $rowData = $data[$rowNum]
    foreach ($col in $rowData.Keys) {
        $addr = "$col$rowNum"
        $ws.Range($addr).Value2 = $rowData[$col]
    }
}

Write-Output "Updated $($data.Count) rows across columns B,C,D,E,F,I,J,K,L,M"